$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 46
$ws1.Range("F3").Value = 326
$ws1.Range("F4").Value = 258
$ws1.Range("F5").Value = 3031
$ws1.Range("F6").Value = 2034
$ws1.Range("F8").Value = 141
$ws1.Range("F9").Value = 1134
$ws1.Range("F10").Value = 208
$ws1.Range("F11").Value = 767
$ws1.Range("F12").Value = 67

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 28

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 46
$ws4.Range("F3").Value = 326
$ws4.Range("F4").Value = 258
$ws4.Range("F5").Value = 3031
$ws4.Range("F6").Value = 2034
$ws4.Range("F8").Value = 28
$ws4.Range("F9").Value = 141
$ws4.Range("F10").Value = 1134
$ws4.Range("F11").Value = 208
$ws4.Range("F12").Value = 768
$ws4.Range("F13").Value = 67

$wb.Save()
